# Update Name of Algo
# Applies refreshed imputation values (RandomForest re-run) to columns A and B
# of the active worksheet, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -22.04730000000001
$ws.Range("B4").Value = 5.292400000000002
$ws.Range("A6").Value = -22.81120000000002
$ws.Range("A7").Value = -19.81849999999999
$ws.Range("A8").Value = -22.32210000000002
$ws.Range("B8").Value = 5.078300000000004
$ws.Range("B9").Value = 6.161300000000001
$ws.Range("B12").Value = 4.795799999999996
$ws.Range("A16").Value = -21.92420000000001
$ws.Range("B17").Value = 5.524
$ws.Range("B18").Value = 6.330199999999996
$ws.Range("B19").Value = 9.208399999999994
$ws.Range("A20").Value = -22.84010000000002
$ws.Range("B20").Value = 5.246299999999996
$ws.Range("A21").Value = -20.00059999999998
$ws.Range("B26").Value = 4.302500000000004
$ws.Range("A28").Value = -22.01260000000001
$ws.Range("A29").Value = -21.32419999999997
$ws.Range("A30").Value = -21.81130000000001
$ws.Range("B31").Value = 4.642799999999998
$ws.Range("A32").Value = -21.1816
$ws.Range("B39").Value = 9.343800000000005
$ws.Range("A40").Value = -21.78509999999998
$ws.Range("B40").Value = 6.098499999999999
$ws.Range("B41").Value = 9.355799999999988
$ws.Range("B42").Value = 9.305699999999991
$ws.Range("B43").Value = 6.028100000000002
$ws.Range("A46").Value = -21.82370000000001
$ws.Range("B47").Value = 5.650700000000003
$ws.Range("B48").Value = 5.430700000000003
$ws.Range("A51").Value = -21.58599999999998
$ws.Range("A52").Value = -21.954
$ws.Range("B54").Value = 4.873900000000003
$ws.Range("A57").Value = -21.97910000000002
$ws.Range("A59").Value = -22.22240000000001
$ws.Range("A62").Value = -22.40130000000002
$ws.Range("B62").Value = 4.989900000000002
$ws.Range("B63").Value = 4.830099999999998
$ws.Range("B64").Value = 5.406799999999999
$ws.Range("A66").Value = -21.4112
$ws.Range("A73").Value = -20.15969999999999
$ws.Range("A74").Value = -21.58839999999998
$ws.Range("B76").Value = 5.769899999999999
$ws.Range("A77").Value = -20.09769999999999
$ws.Range("B81").Value = 5.419600000000003
$ws.Range("B84").Value = 5.745700000000002
$ws.Range("B89").Value = 5.489899999999998
$ws.Range("A92").Value = -21.6268
$ws.Range("B94").Value = 4.723099999999993
$ws.Range("A100").Value = -22.08150000000001
